$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95 (shifts existing rows 95-134 down to 96-135)
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new weekly record
$ws.Cells.Item(95, 1).Value = 4
$ws.Cells.Item(95, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(95, 3).Value = "Los Lagos"
$ws.Cells.Item(95, 4).Value = 45006
$ws.Cells.Item(95, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(95, 5).Value = 10
$ws.Cells.Item(95, 6).Value = 100112031
$ws.Cells.Item(95, 7).Value = "Poroto verde"
$ws.Cells.Item(95, 8).Value = "Magnum"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 45
$ws.Cells.Item(95, 11).Value = 32000
$ws.Cells.Item(95, 12).Value = 32000
$ws.Cells.Item(95, 13).Value = 32000
$ws.Cells.Item(95, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(95, 15).Value = "Región Metropolitana"
$ws.Cells.Item(95, 16).Value = 1280
$ws.Cells.Item(95, 17).Value = 25
$ws.Cells.Item(95, 18).Value = "Hortaliza"
